$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 02:46"

# --- Update row 6 stats ---
$ws.Range("B6").Value = 27069
$ws.Range("C6").Value = 7686
$ws.Range("E6").Value = 26553
$ws.Range("G6").Value = 84
$ws.Range("H6").Value = 340

# --- Update row 53 stats ---
$ws.Range("B53").Value = 274
$ws.Range("C53").Value = 18
$ws.Range("D53").Value = 100
$ws.Range("E53").Value = 170

# --- Reorder Ghana / Cuba ahead of Puerto Rico (rows 116-118) and update Puerto Rico stats ---
# Row 116: Ghana (was Puerto Rico)
$ws.Range("A116").Value = "Ghana"
$ws.Range("C116").Value = 5
$ws.Range("E116").Value = 20
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 1

# Row 117: Cuba (was Ghana)
$ws.Range("A117").Value = "Cuba"
$ws.Range("C117").Value = 0
$ws.Range("E117").Value = 20
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 1

# Row 118: Puerto Rico (was Cuba), with updated stats
$ws.Range("A118").Value = "Puerto Rico"
$ws.Range("C118").Value = 7
$ws.Range("E118").Value = 20
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 1
